$wb = $excel.ActiveWorkbook

# --- Rename the "Include from LOINC" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item("Include from LOINC")
$wsInclude.Name = "Include #0"

# --- Update values on the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version 0.1.0 -> 0.1.1
$wsMeta.Range("B3").Value = "0.1.1"

# Date updated
$wsMeta.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new "Jurisdiction" row right after "Contact" (row 10), before "Description" (old row 11)
$wsMeta.Rows.Item(11).Insert()

# Match formatting of the surrounding data rows for the newly inserted row
$wsMeta.Range("A12:B12").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
